$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 values get entered first (matches shared-string insertion order seen
# in the target file: "Slow Log File Writes" = idx 32, "Slow Data File Reads"
# = idx 33, then the shared URL = idx 34).
$ws.Range("A16").Value = 12
$ws.Range("B16").Value = 50
$ws.Range("C16").Value = "Server Performance"
$ws.Range("D16").Value = "Slow Log File Writes"

$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 50
$ws.Range("C15").Value = "Server Performance"
$ws.Range("D15").Value = "Slow Data File Reads"

$ws.Range("E15").Value = "http://BrentOzar.com/go/slow"
$ws.Hyperlinks.Add($ws.Range("E15"), "http://BrentOzar.com/go/slow")
$ws.Range("E15").Style = "Hyperlink"

$ws.Range("E16").Value = "http://BrentOzar.com/go/slow"
$ws.Hyperlinks.Add($ws.Range("E16"), "http://BrentOzar.com/go/slow")
$ws.Range("E16").Style = "Hyperlink"

$ws.Range("A15").Select()
